$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 65799
$ws.Range("J3").Value = 65799
$ws.Range("L3").Value = 65799
$ws.Range("N3").Value = -66027
# Row 33
$ws.Range("H33").Value = 439.8
$ws.Range("I33").Value = 439.8
$ws.Range("K33").Value = 439.8
$ws.Range("M33").Value = -210.8
# Row 40
$ws.Range("H40").Value = 4821.5
$ws.Range("J40").Value = 5385.9
$ws.Range("L40").Value = 5385.9
$ws.Range("N40").Value = -5735.9
# Row 74
$ws.Range("H74").Value = 6368.5
$ws.Range("I74").Value = 5074.6665
$ws.Range("K74").Value = 5074.6665
$ws.Range("M74").Value = -4138.6665
# Row 77
$ws.Range("H77").Value = 6368.5
$ws.Range("I77").Value = 5074.6665
$ws.Range("K77").Value = 25373.3325
$ws.Range("M77").Value = -20693.3325
# Row 88
$ws.Range("H88").Value = 4752.923
$ws.Range("J88").Value = 5024
$ws.Range("L88").Value = 5024
$ws.Range("N88").Value = -5836
# Row 91
$ws.Range("H91").Value = 4752.923
$ws.Range("J91").Value = 5024
$ws.Range("L91").Value = 5024
$ws.Range("N91").Value = -7832
# Row 102
$ws.Range("H102").Value = 65799
$ws.Range("J102").Value = 65799
$ws.Range("L102").Value = 65799
$ws.Range("N102").Value = -72289

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 44
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
# Row 55
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
# Row 80
$ws.Range("H80").Value = 33992
$ws.Range("J80").Value = 33992.25
$ws.Range("L80").Value = 33992.25
$ws.Range("N80").Value = -35988.25
# Row 83
$ws.Range("H83").Value = 33992
$ws.Range("J83").Value = 33992.25
$ws.Range("L83").Value = 101976.75
$ws.Range("N83").Value = -111960.75
# Row 88
$ws.Range("H88").Value = 2400
$ws.Range("I88").Value = 1431.3334
$ws.Range("J88").Value = 3126.5
$ws.Range("K88").Value = 1431.3334
$ws.Range("L88").Value = 3126.5
$ws.Range("M88").Value = -1025.3334
$ws.Range("N88").Value = -3938.5
# Row 91
$ws.Range("H91").Value = 2400
$ws.Range("I91").Value = 1431.3334
$ws.Range("J91").Value = 3126.5
$ws.Range("K91").Value = 1431.3334
$ws.Range("L91").Value = 3126.5
$ws.Range("M91").Value = -27.33339999999998
$ws.Range("N91").Value = -5934.5
# Row 101
$ws.Range("H101").Value = 28533.666
$ws.Range("J101").Value = 28533.666
$ws.Range("L101").Value = 28533.666
$ws.Range("N101").Value = -35023.666

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 15638.25
$ws.Range("J7").Value = 31177
$ws.Range("L7").Value = 31177
$ws.Range("N7").Value = -31403
# Row 82
$ws.Range("H82").Value = 42542
$ws.Range("J82").Value = 69998.336
$ws.Range("L82").Value = 69998.336
$ws.Range("N82").Value = -70764.336
# Row 85
$ws.Range("H85").Value = 42542
$ws.Range("J85").Value = 69998.336
$ws.Range("L85").Value = 69998.336
$ws.Range("N85").Value = -72650.336
# Row 103
$ws.Range("H103").Value = 17999.666
$ws.Range("J103").Value = 17999.666
$ws.Range("L103").Value = 17999.666
$ws.Range("N103").Value = -20343.666
# Row 134
$ws.Range("H134").Value = 3940.7144
$ws.Range("I134").Value = 2994.4
$ws.Range("J134").Value = 6306.5
$ws.Range("K134").Value = 8983.200000000001
$ws.Range("L134").Value = 18919.5
$ws.Range("M134").Value = -6448.200000000001
$ws.Range("N134").Value = -23989.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 60
$ws.Range("H60").Value = 39333
$ws.Range("I60").Value = 18666.666
$ws.Range("K60").Value = 18666.666
$ws.Range("M60").Value = -18155.666

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 44
$ws.Range("H44").Value = 808.1667
$ws.Range("I44").Value = 808.1667
$ws.Range("K44").Value = 2424.5001
$ws.Range("M44").Value = -2026.5001
# Row 76
$ws.Range("H76").Value = 113
$ws.Range("I76").Value = 113
$ws.Range("K76").Value = 339
$ws.Range("M76").Value = 44
# Row 79
$ws.Range("H79").Value = 113
$ws.Range("I79").Value = 113
$ws.Range("K79").Value = 339
$ws.Range("M79").Value = 987
# Row 108
$ws.Range("H108").Value = 636.25
$ws.Range("I108").Value = 636.25
$ws.Range("K108").Value = 1908.75
$ws.Range("M108").Value = 971.25
# Row 115
$ws.Range("H115").Value = 733.3333
$ws.Range("I115").Value = 733.3333
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 2199.9999
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -1024.9999
$ws.Range("N115").ClearContents()
# Row 129
$ws.Range("H129").Value = 2482.4443
$ws.Range("I129").Value = 345
$ws.Range("J129").Value = 3093.1428
$ws.Range("K129").Value = 1035
$ws.Range("L129").Value = 9279.428400000001
$ws.Range("M129").Value = 3965
$ws.Range("N129").Value = -19279.4284
# Row 132
$ws.Range("H132").Value = 1702
$ws.Range("I132").Value = 1702
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 15318
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -12788
$ws.Range("N132").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 21
$ws.Range("H21").Value = 8699999
$ws.Range("J21").Value = 49999
$ws.Range("L21").Value = 49999
$ws.Range("N21").Value = -50345
# Row 30
$ws.Range("H30").Value = 8699999
$ws.Range("J30").Value = 49999
$ws.Range("L30").Value = 49999
$ws.Range("N30").Value = -50209
# Row 122
$ws.Range("H122").Value = 2840.4285
$ws.Range("I122").Value = 2397.353
$ws.Range("J122").Value = 4723.5
$ws.Range("K122").Value = 7192.059
$ws.Range("L122").Value = 14170.5
$ws.Range("M122").Value = -4742.059
$ws.Range("N122").Value = -19070.5
# Row 132
$ws.Range("H132").Value = 2235.3684
$ws.Range("I132").Value = 2380.25
$ws.Range("K132").Value = 7140.75
$ws.Range("M132").Value = -4610.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 8
$ws.Range("H8").Value = 120000
$ws.Range("J8").Value = 120000
$ws.Range("L8").Value = 120000
$ws.Range("N8").Value = -120280
# Row 22
$ws.Range("H22").Value = 1868.375
$ws.Range("I22").Value = 1699.4
$ws.Range("K22").Value = 1699.4
$ws.Range("M22").Value = -1404.4
# Row 27
$ws.Range("H27").Value = 1868.375
$ws.Range("I27").Value = 1699.4
$ws.Range("K27").Value = 1699.4
$ws.Range("M27").Value = -1592.4
# Row 32
$ws.Range("H32").Value = 4800
$ws.Range("I32").Value = 2766.6667
$ws.Range("J32").Value = 10900
$ws.Range("K32").Value = 2766.6667
$ws.Range("L32").Value = 10900
$ws.Range("M32").Value = -2449.6667
$ws.Range("N32").Value = -11534
# Row 100
$ws.Range("H100").Value = 2715.6667
$ws.Range("I100").Value = 2573.5
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 2573.5
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -2032.5
$ws.Range("N100").Value = -4082
# Row 118
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
# Row 136
$ws.Range("H136").Value = 1549
$ws.Range("I136").Value = 1747
$ws.Range("J136").Value = 1450
$ws.Range("K136").Value = 5241
$ws.Range("L136").Value = 4350
$ws.Range("M136").Value = -2691
$ws.Range("N136").Value = -9450

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Range("H45").Value = 39999.5
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 39999.5
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 39999.5
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -40981.5
# Row 81
$ws.Range("H81").Value = 2250
$ws.Range("I81").Value = 2250
$ws.Range("K81").Value = 4500
$ws.Range("M81").Value = -3439
# Row 84
$ws.Range("H84").Value = 2250
$ws.Range("I84").Value = 2250
$ws.Range("K84").Value = 22500
$ws.Range("M84").Value = -17196
# Row 88
$ws.Range("H88").Value = 39000
$ws.Range("I88").Value = 18000
$ws.Range("J88").Value = 60000
$ws.Range("K88").Value = 18000
$ws.Range("L88").Value = 60000
$ws.Range("M88").Value = -17594
$ws.Range("N88").Value = -60812
# Row 91
$ws.Range("H91").Value = 39000
$ws.Range("I91").Value = 18000
$ws.Range("J91").Value = 60000
$ws.Range("K91").Value = 18000
$ws.Range("L91").Value = 60000
$ws.Range("M91").Value = -16596
$ws.Range("N91").Value = -62808
# Row 136
$ws.Range("H136").Value = 9049
$ws.Range("I136").Value = 9049
$ws.Range("K136").Value = 27147
$ws.Range("M136").Value = -24597
